# Add a new row of data (row 8) to Sheet1, describing the
# "Binary Tree Level Order Traversal" problem, matching the
# existing table style used by rows 6 and 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new row's values.
$ws.Range("A8").Value = "Binary Tree Level Order Traversal"
$ws.Range("B8").Value = "Return BT traversed in list of list"
$ws.Range("C8").Value = "Use queue to iteratively BFS traverse. Keep count of levels and add values to list accordingly."

# Copy formatting from the row above (row 7) so the new row matches
# the existing "Name" / "Description" style used by the table.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Update the active selection to reflect the new last-edited cell.
$ws.Range("C8").Select()
